{"js": "// Update the equation values in the first table of the document body.\n// The table holds 20 rows x 5 columns of \"a OP b = c\" arithmetic strings;\n// each cell is updated in place to the new equation text. Assigning a\n// full 2-D array to Table.values rewrites each cell's text while Word\n// keeps the cell's existing paragraph/run formatting (font, size, etc.)\n// untouched, since only the run text is modified.\n\nconst NEW_VALUES = [\n  [\"94-44=50\", \"61-58=3\", \"14-11=3\", \"39+60=99\", \"3+20=23\"],\n  [\"9+39=48\", \"20+18=38\", \"6+67=73\", \"9+52=61\", \"71-34=37\"],\n  [\"49+44=93\", \"23+45=68\", \"61-28=33\", \"32-0=32\", \"23+44=67\"],\n  [\"57-38=19\", \"76+3=79\", \"55+37=92\", \"67+21=88\", \"47+19=66\"],\n  [\"56-22=34\", \"88-24=64\", \"71+20=91\", \"17+77=94\", \"92-63=29\"],\n  [\"20+74=94\", \"36-29=7\", \"49+23=72\", \"56+27=83\", \"57+33=90\"],\n  [\"36+57=93\", \"59+36=95\", \"82-15=67\", \"67+12=79\", \"84-65=19\"],\n  [\"95-2=93\", \"13+57=70\", \"94-86=8\", \"62+36=98\", \"34-1=33\"],\n  [\"65+6=71\", \"86-56=30\", \"14-7=7\", \"59-12=47\", \"4+16=20\"],\n  [\"96-26=70\", \"75-37=38\", \"62-41=21\", \"54+38=92\", \"41-23=18\"],\n  [\"45+43=88\", \"26+32=58\", \"82-8=74\", \"25+9=34\", \"64-52=12\"],\n  [\"70-43=27\", \"62+22=84\", \"70-57=13\", \"8+41=49\", \"11+37=48\"],\n  [\"36+11=47\", \"13+51=64\", \"17+50=67\", \"19+78=97\", \"90-52=38\"],\n  [\"78-12=66\", \"25+43=68\", \"29-1=28\", \"53-28=25\", \"29-18=11\"],\n  [\"39+7=46\", \"49-1=48\", \"9+57=66\", \"54+15=69\", \"35+47=82\"],\n  [\"34+22=56\", \"1+96=97\", \"21+74=95\", \"28+51=79\", \"23+28=51\"],\n  [\"38-3=35\", \"80-2=78\", \"15+56=71\", \"54+45=99\", \"1+45=46\"],\n  [\"51-20=31\", \"50-16=34\", \"61-31=30\", \"98-20=78\", \"57-13=44\"],\n  [\"70+15=85\", \"28+49=77\", \"29+9=38\", \"98-82=16\", \"95-94=1\"],\n  [\"60-2=58\", \"13-1=12\", \"16+40=56\", \"45-22=23\", \"49+2=51\"]\n];\n\nconst tables = context.document.body.tables;\ntables.load(\"items\");\nawait context.sync();\n\nif (tables.items.length === 0) {\n  throw new Error(\"Expected at least one table in the document body.\");\n}\n\nconst table = tables.items[0];\ntable.values = NEW_VALUES;\nawait context.sync();\n", "ps1": "# Update the 20x5 grid of arithmetic equations in the first table of the\n# document body. Each cell is addressed by its (row, column) position so\n# the update is correct even though some original equation strings repeat\n# (e.g. \"37+12=49\" appears twice with different replacements).\n$d = $word.ActiveDocument\n$t = $d.Tables.Item(1)\n\n$newValues = @(\n    @(\"94-44=50\", \"61-58=3\", \"14-11=3\", \"39+60=99\", \"3+20=23\"),\n    @(\"9+39=48\", \"20+18=38\", \"6+67=73\", \"9+52=61\", \"71-34=37\"),\n    @(\"49+44=93\", \"23+45=68\", \"61-28=33\", \"32-0=32\", \"23+44=67\"),\n    @(\"57-38=19\", \"76+3=79\", \"55+37=92\", \"67+21=88\", \"47+19=66\"),\n    @(\"56-22=34\", \"88-24=64\", \"71+20=91\", \"17+77=94\", \"92-63=29\"),\n    @(\"20+74=94\", \"36-29=7\", \"49+23=72\", \"56+27=83\", \"57+33=90\"),\n    @(\"36+57=93\", \"59+36=95\", \"82-15=67\", \"67+12=79\", \"84-65=19\"),\n    @(\"95-2=93\", \"13+57=70\", \"94-86=8\", \"62+36=98\", \"34-1=33\"),\n    @(\"65+6=71\", \"86-56=30\", \"14-7=7\", \"59-12=47\", \"4+16=20\"),\n    @(\"96-26=70\", \"75-37=38\", \"62-41=21\", \"54+38=92\", \"41-23=18\"),\n    @(\"45+43=88\", \"26+32=58\", \"82-8=74\", \"25+9=34\", \"64-52=12\"),\n    @(\"70-43=27\", \"62+22=84\", \"70-57=13\", \"8+41=49\", \"11+37=48\"),\n    @(\"36+11=47\", \"13+51=64\", \"17+50=67\", \"19+78=97\", \"90-52=38\"),\n    @(\"78-12=66\", \"25+43=68\", \"29-1=28\", \"53-28=25\", \"29-18=11\"),\n    @(\"39+7=46\", \"49-1=48\", \"9+57=66\", \"54+15=69\", \"35+47=82\"),\n    @(\"34+22=56\", \"1+96=97\", \"21+74=95\", \"28+51=79\", \"23+28=51\"),\n    @(\"38-3=35\", \"80-2=78\", \"15+56=71\", \"54+45=99\", \"1+45=46\"),\n    @(\"51-20=31\", \"50-16=34\", \"61-31=30\", \"98-20=78\", \"57-13=44\"),\n    @(\"70+15=85\", \"28+49=77\", \"29+9=38\", \"98-82=16\", \"95-94=1\"),\n    @(\"60-2=58\", \"13-1=12\", \"16+40=56\", \"45-22=23\", \"49+2=51\"),\n)\n\n$rows = $t.Rows.Count\n$cols = $t.Columns.Count\nfor ($r = 1; $r -le $rows; $r++) {\n    for ($c = 1; $c -le $cols; $c++) {\n        $cell = $t.Cell($r, $c)\n        $cell.Range.Text = $newValues[$r - 1][$c - 1]\n    }\n}\n"}
